$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 352
$ws.Cells.Item(352, 1).Value = 11
$ws.Cells.Item(352, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(352, 3).Value = 'Bíobío'
$ws.Cells.Item(352, 4).Value = 44516
$ws.Cells.Item(352, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(352, 5).Value = 8
$ws.Cells.Item(352, 6).Value = 'Fruta'
$ws.Cells.Item(352, 7).Value = 100108
$ws.Cells.Item(352, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(352, 9).Value = 100108006
$ws.Cells.Item(352, 10).Value = 'Plátano'
$ws.Cells.Item(352, 11).Value = 'Sin especificar'
$ws.Cells.Item(352, 12).Value = 'Maduro'
$ws.Cells.Item(352, 13).Value = 200
$ws.Cells.Item(352, 14).Value = 11000
$ws.Cells.Item(352, 15).Value = 11000
$ws.Cells.Item(352, 16).Value = 11000
$ws.Cells.Item(352, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(352, 18).Value = 'Ecuador'
$ws.Cells.Item(352, 19).Value = 550
$ws.Cells.Item(352, 20).Value = 20

# Row 353
$ws.Cells.Item(353, 1).Value = 11
$ws.Cells.Item(353, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(353, 3).Value = 'Bíobío'
$ws.Cells.Item(353, 4).Value = 44516
$ws.Cells.Item(353, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(353, 5).Value = 8
$ws.Cells.Item(353, 6).Value = 'Fruta'
$ws.Cells.Item(353, 7).Value = 100108
$ws.Cells.Item(353, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(353, 9).Value = 100108006
$ws.Cells.Item(353, 10).Value = 'Plátano'
$ws.Cells.Item(353, 11).Value = 'Sin especificar'
$ws.Cells.Item(353, 12).Value = 'Pintón'
$ws.Cells.Item(353, 13).Value = 400
$ws.Cells.Item(353, 14).Value = 13000
$ws.Cells.Item(353, 15).Value = 13000
$ws.Cells.Item(353, 16).Value = 13000
$ws.Cells.Item(353, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(353, 18).Value = 'Ecuador'
$ws.Cells.Item(353, 19).Value = 650
$ws.Cells.Item(353, 20).Value = 20

# Row 354
$ws.Cells.Item(354, 1).Value = 11
$ws.Cells.Item(354, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(354, 3).Value = 'Bíobío'
$ws.Cells.Item(354, 4).Value = 44516
$ws.Cells.Item(354, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(354, 5).Value = 8
$ws.Cells.Item(354, 6).Value = 'Fruta'
$ws.Cells.Item(354, 7).Value = 100108
$ws.Cells.Item(354, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(354, 9).Value = 100108006
$ws.Cells.Item(354, 10).Value = 'Plátano'
$ws.Cells.Item(354, 11).Value = 'Sin especificar'
$ws.Cells.Item(354, 12).Value = 'Primera Pintón'
$ws.Cells.Item(354, 13).Value = 300
$ws.Cells.Item(354, 14).Value = 15000
$ws.Cells.Item(354, 15).Value = 15000
$ws.Cells.Item(354, 16).Value = 15000
$ws.Cells.Item(354, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(354, 18).Value = 'Ecuador'
$ws.Cells.Item(354, 19).Value = 750
$ws.Cells.Item(354, 20).Value = 20

# Row 355
$ws.Cells.Item(355, 1).Value = 11
$ws.Cells.Item(355, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(355, 3).Value = 'Bíobío'
$ws.Cells.Item(355, 4).Value = 44295
$ws.Cells.Item(355, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(355, 5).Value = 8
$ws.Cells.Item(355, 6).Value = 'Fruta'
$ws.Cells.Item(355, 7).Value = 100108
$ws.Cells.Item(355, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(355, 9).Value = 100108006
$ws.Cells.Item(355, 10).Value = 'Plátano'
$ws.Cells.Item(355, 11).Value = 'Sin especificar'
$ws.Cells.Item(355, 12).Value = 'Maduro'
$ws.Cells.Item(355, 13).Value = 100
$ws.Cells.Item(355, 14).Value = 11000
$ws.Cells.Item(355, 15).Value = 11000
$ws.Cells.Item(355, 16).Value = 11000
$ws.Cells.Item(355, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(355, 18).Value = 'Ecuador'
$ws.Cells.Item(355, 19).Value = 550
$ws.Cells.Item(355, 20).Value = 20

# Row 356
$ws.Cells.Item(356, 1).Value = 11
$ws.Cells.Item(356, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(356, 3).Value = 'Bíobío'
$ws.Cells.Item(356, 4).Value = 44295
$ws.Cells.Item(356, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(356, 5).Value = 8
$ws.Cells.Item(356, 6).Value = 'Fruta'
$ws.Cells.Item(356, 7).Value = 100108
$ws.Cells.Item(356, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(356, 9).Value = 100108006
$ws.Cells.Item(356, 10).Value = 'Plátano'
$ws.Cells.Item(356, 11).Value = 'Sin especificar'
$ws.Cells.Item(356, 12).Value = 'Pintón'
$ws.Cells.Item(356, 13).Value = 300
$ws.Cells.Item(356, 14).Value = 13000
$ws.Cells.Item(356, 15).Value = 13000
$ws.Cells.Item(356, 16).Value = 13000
$ws.Cells.Item(356, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(356, 18).Value = 'Ecuador'
$ws.Cells.Item(356, 19).Value = 650
$ws.Cells.Item(356, 20).Value = 20

# Row 357
$ws.Cells.Item(357, 1).Value = 11
$ws.Cells.Item(357, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(357, 3).Value = 'Bíobío'
$ws.Cells.Item(357, 4).Value = 44295
$ws.Cells.Item(357, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(357, 5).Value = 8
$ws.Cells.Item(357, 6).Value = 'Fruta'
$ws.Cells.Item(357, 7).Value = 100108
$ws.Cells.Item(357, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(357, 9).Value = 100108006
$ws.Cells.Item(357, 10).Value = 'Plátano'
$ws.Cells.Item(357, 11).Value = 'Sin especificar'
$ws.Cells.Item(357, 12).Value = 'Primera Pintón'
$ws.Cells.Item(357, 13).Value = 300
$ws.Cells.Item(357, 14).Value = 15000
$ws.Cells.Item(357, 15).Value = 15000
$ws.Cells.Item(357, 16).Value = 15000
$ws.Cells.Item(357, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(357, 18).Value = 'Ecuador'
$ws.Cells.Item(357, 19).Value = 750
$ws.Cells.Item(357, 20).Value = 20

# Row 358
$ws.Cells.Item(358, 1).Value = 11
$ws.Cells.Item(358, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(358, 3).Value = 'Bíobío'
$ws.Cells.Item(358, 4).Value = 44217
$ws.Cells.Item(358, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(358, 5).Value = 8
$ws.Cells.Item(358, 6).Value = 'Fruta'
$ws.Cells.Item(358, 7).Value = 100108
$ws.Cells.Item(358, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(358, 9).Value = 100108006
$ws.Cells.Item(358, 10).Value = 'Plátano'
$ws.Cells.Item(358, 11).Value = 'Sin especificar'
$ws.Cells.Item(358, 12).Value = 'Maduro'
$ws.Cells.Item(358, 13).Value = 200
$ws.Cells.Item(358, 14).Value = 9000
$ws.Cells.Item(358, 15).Value = 9000
$ws.Cells.Item(358, 16).Value = 9000
$ws.Cells.Item(358, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(358, 18).Value = 'Ecuador'
$ws.Cells.Item(358, 19).Value = 450
$ws.Cells.Item(358, 20).Value = 20

# Row 359
$ws.Cells.Item(359, 1).Value = 11
$ws.Cells.Item(359, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(359, 3).Value = 'Bíobío'
$ws.Cells.Item(359, 4).Value = 44217
$ws.Cells.Item(359, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(359, 5).Value = 8
$ws.Cells.Item(359, 6).Value = 'Fruta'
$ws.Cells.Item(359, 7).Value = 100108
$ws.Cells.Item(359, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(359, 9).Value = 100108006
$ws.Cells.Item(359, 10).Value = 'Plátano'
$ws.Cells.Item(359, 11).Value = 'Sin especificar'
$ws.Cells.Item(359, 12).Value = 'Pintón'
$ws.Cells.Item(359, 13).Value = 400
$ws.Cells.Item(359, 14).Value = 11000
$ws.Cells.Item(359, 15).Value = 11000
$ws.Cells.Item(359, 16).Value = 11000
$ws.Cells.Item(359, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(359, 18).Value = 'Ecuador'
$ws.Cells.Item(359, 19).Value = 550
$ws.Cells.Item(359, 20).Value = 20

# Row 360
$ws.Cells.Item(360, 1).Value = 11
$ws.Cells.Item(360, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(360, 3).Value = 'Bíobío'
$ws.Cells.Item(360, 4).Value = 44217
$ws.Cells.Item(360, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(360, 5).Value = 8
$ws.Cells.Item(360, 6).Value = 'Fruta'
$ws.Cells.Item(360, 7).Value = 100108
$ws.Cells.Item(360, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(360, 9).Value = 100108006
$ws.Cells.Item(360, 10).Value = 'Plátano'
$ws.Cells.Item(360, 11).Value = 'Sin especificar'
$ws.Cells.Item(360, 12).Value = 'Primera Pintón'
$ws.Cells.Item(360, 13).Value = 400
$ws.Cells.Item(360, 14).Value = 13000
$ws.Cells.Item(360, 15).Value = 13000
$ws.Cells.Item(360, 16).Value = 13000
$ws.Cells.Item(360, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(360, 18).Value = 'Ecuador'
$ws.Cells.Item(360, 19).Value = 650
$ws.Cells.Item(360, 20).Value = 20

# Row 361
$ws.Cells.Item(361, 1).Value = 11
$ws.Cells.Item(361, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(361, 3).Value = 'Bíobío'
$ws.Cells.Item(361, 4).Value = 44509
$ws.Cells.Item(361, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(361, 5).Value = 8
$ws.Cells.Item(361, 6).Value = 'Fruta'
$ws.Cells.Item(361, 7).Value = 100108
$ws.Cells.Item(361, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(361, 9).Value = 100108006
$ws.Cells.Item(361, 10).Value = 'Plátano'
$ws.Cells.Item(361, 11).Value = 'Sin especificar'
$ws.Cells.Item(361, 12).Value = 'Primera Pintón'
$ws.Cells.Item(361, 13).Value = 1100
$ws.Cells.Item(361, 14).Value = 17000
$ws.Cells.Item(361, 15).Value = 18000
$ws.Cells.Item(361, 16).Value = 17545
$ws.Cells.Item(361, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(361, 18).Value = 'Ecuador'
$ws.Cells.Item(361, 19).Value = 877
$ws.Cells.Item(361, 20).Value = 20

# Row 362
$ws.Cells.Item(362, 1).Value = 11
$ws.Cells.Item(362, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(362, 3).Value = 'Bíobío'
$ws.Cells.Item(362, 4).Value = 44421
$ws.Cells.Item(362, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(362, 5).Value = 8
$ws.Cells.Item(362, 6).Value = 'Fruta'
$ws.Cells.Item(362, 7).Value = 100108
$ws.Cells.Item(362, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(362, 9).Value = 100108006
$ws.Cells.Item(362, 10).Value = 'Plátano'
$ws.Cells.Item(362, 11).Value = 'Sin especificar'
$ws.Cells.Item(362, 12).Value = 'Maduro'
$ws.Cells.Item(362, 13).Value = 100
$ws.Cells.Item(362, 14).Value = 10000
$ws.Cells.Item(362, 15).Value = 10000
$ws.Cells.Item(362, 16).Value = 10000
$ws.Cells.Item(362, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(362, 18).Value = 'Ecuador'
$ws.Cells.Item(362, 19).Value = 500
$ws.Cells.Item(362, 20).Value = 20

# Row 363
$ws.Cells.Item(363, 1).Value = 11
$ws.Cells.Item(363, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(363, 3).Value = 'Bíobío'
$ws.Cells.Item(363, 4).Value = 44421
$ws.Cells.Item(363, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(363, 5).Value = 8
$ws.Cells.Item(363, 6).Value = 'Fruta'
$ws.Cells.Item(363, 7).Value = 100108
$ws.Cells.Item(363, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(363, 9).Value = 100108006
$ws.Cells.Item(363, 10).Value = 'Plátano'
$ws.Cells.Item(363, 11).Value = 'Sin especificar'
$ws.Cells.Item(363, 12).Value = 'Pintón'
$ws.Cells.Item(363, 13).Value = 300
$ws.Cells.Item(363, 14).Value = 12000
$ws.Cells.Item(363, 15).Value = 12000
$ws.Cells.Item(363, 16).Value = 12000
$ws.Cells.Item(363, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(363, 18).Value = 'Ecuador'
$ws.Cells.Item(363, 19).Value = 600
$ws.Cells.Item(363, 20).Value = 20

# Row 364
$ws.Cells.Item(364, 1).Value = 11
$ws.Cells.Item(364, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(364, 3).Value = 'Bíobío'
$ws.Cells.Item(364, 4).Value = 44421
$ws.Cells.Item(364, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(364, 5).Value = 8
$ws.Cells.Item(364, 6).Value = 'Fruta'
$ws.Cells.Item(364, 7).Value = 100108
$ws.Cells.Item(364, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(364, 9).Value = 100108006
$ws.Cells.Item(364, 10).Value = 'Plátano'
$ws.Cells.Item(364, 11).Value = 'Sin especificar'
$ws.Cells.Item(364, 12).Value = 'Primera Pintón'
$ws.Cells.Item(364, 13).Value = 300
$ws.Cells.Item(364, 14).Value = 14000
$ws.Cells.Item(364, 15).Value = 14000
$ws.Cells.Item(364, 16).Value = 14000
$ws.Cells.Item(364, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(364, 18).Value = 'Ecuador'
$ws.Cells.Item(364, 19).Value = 700
$ws.Cells.Item(364, 20).Value = 20

# Row 365
$ws.Cells.Item(365, 1).Value = 11
$ws.Cells.Item(365, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(365, 3).Value = 'Bíobío'
$ws.Cells.Item(365, 4).Value = 44383
$ws.Cells.Item(365, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(365, 5).Value = 8
$ws.Cells.Item(365, 6).Value = 'Fruta'
$ws.Cells.Item(365, 7).Value = 100108
$ws.Cells.Item(365, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(365, 9).Value = 100108006
$ws.Cells.Item(365, 10).Value = 'Plátano'
$ws.Cells.Item(365, 11).Value = 'Sin especificar'
$ws.Cells.Item(365, 12).Value = 'Maduro'
$ws.Cells.Item(365, 13).Value = 100
$ws.Cells.Item(365, 14).Value = 8000
$ws.Cells.Item(365, 15).Value = 8000
$ws.Cells.Item(365, 16).Value = 8000
$ws.Cells.Item(365, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(365, 18).Value = 'Ecuador'
$ws.Cells.Item(365, 19).Value = 400
$ws.Cells.Item(365, 20).Value = 20

# Row 366
$ws.Cells.Item(366, 1).Value = 11
$ws.Cells.Item(366, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(366, 3).Value = 'Bíobío'
$ws.Cells.Item(366, 4).Value = 44383
$ws.Cells.Item(366, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(366, 5).Value = 8
$ws.Cells.Item(366, 6).Value = 'Fruta'
$ws.Cells.Item(366, 7).Value = 100108
$ws.Cells.Item(366, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(366, 9).Value = 100108006
$ws.Cells.Item(366, 10).Value = 'Plátano'
$ws.Cells.Item(366, 11).Value = 'Sin especificar'
$ws.Cells.Item(366, 12).Value = 'Pintón'
$ws.Cells.Item(366, 13).Value = 300
$ws.Cells.Item(366, 14).Value = 10000
$ws.Cells.Item(366, 15).Value = 10000
$ws.Cells.Item(366, 16).Value = 10000
$ws.Cells.Item(366, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(366, 18).Value = 'Ecuador'
$ws.Cells.Item(366, 19).Value = 500
$ws.Cells.Item(366, 20).Value = 20

# Row 367
$ws.Cells.Item(367, 1).Value = 11
$ws.Cells.Item(367, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(367, 3).Value = 'Bíobío'
$ws.Cells.Item(367, 4).Value = 44383
$ws.Cells.Item(367, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(367, 5).Value = 8
$ws.Cells.Item(367, 6).Value = 'Fruta'
$ws.Cells.Item(367, 7).Value = 100108
$ws.Cells.Item(367, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(367, 9).Value = 100108006
$ws.Cells.Item(367, 10).Value = 'Plátano'
$ws.Cells.Item(367, 11).Value = 'Sin especificar'
$ws.Cells.Item(367, 12).Value = 'Primera Pintón'
$ws.Cells.Item(367, 13).Value = 300
$ws.Cells.Item(367, 14).Value = 12000
$ws.Cells.Item(367, 15).Value = 12000
$ws.Cells.Item(367, 16).Value = 12000
$ws.Cells.Item(367, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(367, 18).Value = 'Ecuador'
$ws.Cells.Item(367, 19).Value = 600
$ws.Cells.Item(367, 20).Value = 20

# Row 368
$ws.Cells.Item(368, 1).Value = 11
$ws.Cells.Item(368, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(368, 3).Value = 'Bíobío'
$ws.Cells.Item(368, 4).Value = 44244
$ws.Cells.Item(368, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(368, 5).Value = 8
$ws.Cells.Item(368, 6).Value = 'Fruta'
$ws.Cells.Item(368, 7).Value = 100108
$ws.Cells.Item(368, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(368, 9).Value = 100108006
$ws.Cells.Item(368, 10).Value = 'Plátano'
$ws.Cells.Item(368, 11).Value = 'Sin especificar'
$ws.Cells.Item(368, 12).Value = 'Maduro'
$ws.Cells.Item(368, 13).Value = 50
$ws.Cells.Item(368, 14).Value = 10000
$ws.Cells.Item(368, 15).Value = 10000
$ws.Cells.Item(368, 16).Value = 10000
$ws.Cells.Item(368, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(368, 18).Value = 'Ecuador'
$ws.Cells.Item(368, 19).Value = 500
$ws.Cells.Item(368, 20).Value = 20

# Row 369
$ws.Cells.Item(369, 1).Value = 11
$ws.Cells.Item(369, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(369, 3).Value = 'Bíobío'
$ws.Cells.Item(369, 4).Value = 44244
$ws.Cells.Item(369, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(369, 5).Value = 8
$ws.Cells.Item(369, 6).Value = 'Fruta'
$ws.Cells.Item(369, 7).Value = 100108
$ws.Cells.Item(369, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(369, 9).Value = 100108006
$ws.Cells.Item(369, 10).Value = 'Plátano'
$ws.Cells.Item(369, 11).Value = 'Sin especificar'
$ws.Cells.Item(369, 12).Value = 'Pintón'
$ws.Cells.Item(369, 13).Value = 200
$ws.Cells.Item(369, 14).Value = 11000
$ws.Cells.Item(369, 15).Value = 11000
$ws.Cells.Item(369, 16).Value = 11000
$ws.Cells.Item(369, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(369, 18).Value = 'Ecuador'
$ws.Cells.Item(369, 19).Value = 550
$ws.Cells.Item(369, 20).Value = 20

# Row 370
$ws.Cells.Item(370, 1).Value = 11
$ws.Cells.Item(370, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(370, 3).Value = 'Bíobío'
$ws.Cells.Item(370, 4).Value = 44244
$ws.Cells.Item(370, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(370, 5).Value = 8
$ws.Cells.Item(370, 6).Value = 'Fruta'
$ws.Cells.Item(370, 7).Value = 100108
$ws.Cells.Item(370, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(370, 9).Value = 100108006
$ws.Cells.Item(370, 10).Value = 'Plátano'
$ws.Cells.Item(370, 11).Value = 'Sin especificar'
$ws.Cells.Item(370, 12).Value = 'Primera Pintón'
$ws.Cells.Item(370, 13).Value = 200
$ws.Cells.Item(370, 14).Value = 13000
$ws.Cells.Item(370, 15).Value = 13000
$ws.Cells.Item(370, 16).Value = 13000
$ws.Cells.Item(370, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(370, 18).Value = 'Ecuador'
$ws.Cells.Item(370, 19).Value = 650
$ws.Cells.Item(370, 20).Value = 20

# Row 371
$ws.Cells.Item(371, 1).Value = 11
$ws.Cells.Item(371, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(371, 3).Value = 'Bíobío'
$ws.Cells.Item(371, 4).Value = 44273
$ws.Cells.Item(371, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(371, 5).Value = 8
$ws.Cells.Item(371, 6).Value = 'Fruta'
$ws.Cells.Item(371, 7).Value = 100108
$ws.Cells.Item(371, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(371, 9).Value = 100108006
$ws.Cells.Item(371, 10).Value = 'Plátano'
$ws.Cells.Item(371, 11).Value = 'Sin especificar'
$ws.Cells.Item(371, 12).Value = 'Maduro'
$ws.Cells.Item(371, 13).Value = 200
$ws.Cells.Item(371, 14).Value = 10000
$ws.Cells.Item(371, 15).Value = 10000
$ws.Cells.Item(371, 16).Value = 10000
$ws.Cells.Item(371, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(371, 18).Value = 'Ecuador'
$ws.Cells.Item(371, 19).Value = 500
$ws.Cells.Item(371, 20).Value = 20

# Row 372
$ws.Cells.Item(372, 1).Value = 11
$ws.Cells.Item(372, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(372, 3).Value = 'Bíobío'
$ws.Cells.Item(372, 4).Value = 44273
$ws.Cells.Item(372, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(372, 5).Value = 8
$ws.Cells.Item(372, 6).Value = 'Fruta'
$ws.Cells.Item(372, 7).Value = 100108
$ws.Cells.Item(372, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(372, 9).Value = 100108006
$ws.Cells.Item(372, 10).Value = 'Plátano'
$ws.Cells.Item(372, 11).Value = 'Sin especificar'
$ws.Cells.Item(372, 12).Value = 'Pintón'
$ws.Cells.Item(372, 13).Value = 400
$ws.Cells.Item(372, 14).Value = 12000
$ws.Cells.Item(372, 15).Value = 12000
$ws.Cells.Item(372, 16).Value = 12000
$ws.Cells.Item(372, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(372, 18).Value = 'Ecuador'
$ws.Cells.Item(372, 19).Value = 600
$ws.Cells.Item(372, 20).Value = 20

# Row 373
$ws.Cells.Item(373, 1).Value = 11
$ws.Cells.Item(373, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(373, 3).Value = 'Bíobío'
$ws.Cells.Item(373, 4).Value = 44273
$ws.Cells.Item(373, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(373, 5).Value = 8
$ws.Cells.Item(373, 6).Value = 'Fruta'
$ws.Cells.Item(373, 7).Value = 100108
$ws.Cells.Item(373, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(373, 9).Value = 100108006
$ws.Cells.Item(373, 10).Value = 'Plátano'
$ws.Cells.Item(373, 11).Value = 'Sin especificar'
$ws.Cells.Item(373, 12).Value = 'Primera Pintón'
$ws.Cells.Item(373, 13).Value = 400
$ws.Cells.Item(373, 14).Value = 14000
$ws.Cells.Item(373, 15).Value = 14000
$ws.Cells.Item(373, 16).Value = 14000
$ws.Cells.Item(373, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(373, 18).Value = 'Ecuador'
$ws.Cells.Item(373, 19).Value = 700
$ws.Cells.Item(373, 20).Value = 20
